# Apply crypto price/volume/name updates to sheet1 of the workbook.
# Cells whose new text looks like a plain number need NumberFormat "@"
# forced first, otherwise Excel silently stores them as floating point
# numbers (introducing rounding noise) instead of the literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$ref, [string]$val) {
    $cell = $ws.Range($ref)
    if ($val -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number - force text storage so Excel keeps the
        # exact literal instead of silently converting it to a float.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

Set-TextCell 'D2' '37.044.42'
Set-TextCell 'E2' '  -0.29%  '
Set-TextCell 'D3' '2.003.34'
Set-TextCell 'E3' '  -0.72%  '
Set-TextCell 'E4' '  +0.02%  '
Set-TextCell 'D5' '257.36'
Set-TextCell 'E5' '  +4.50%  '
Set-TextCell 'D6' '0.617'
Set-TextCell 'E6' '  -1.42%  '
Set-TextCell 'E7' '  +0.11%  '
Set-TextCell 'D8' '55.81'
Set-TextCell 'E8' '  -6.75%  '
Set-TextCell 'D9' '0.376'
Set-TextCell 'E9' '  -3.64%  '
Set-TextCell 'D10' '0.0765'
Set-TextCell 'E10' '  -4.92%  '
Set-TextCell 'E11' '  -2.53%  '
Set-TextCell 'D12' '2.300.24'
Set-TextCell 'D13' '14.17'
Set-TextCell 'E13' '  -5.27%  '
Set-TextCell 'D14' '21.28'
Set-TextCell 'E14' '  -2.76%  '
Set-TextCell 'E15' '  -6.39%  '
Set-TextCell 'D16' '5.17'
Set-TextCell 'E16' '  -4.63%  '
Set-TextCell 'D17' '2.001.87'
Set-TextCell 'E17' '  -0.82%  '
Set-TextCell 'D18' '36.998.16'
Set-TextCell 'E18' '  -0.47%  '
Set-TextCell 'D19' '70.79'
Set-TextCell 'E19' '  +0.79%  '
Set-TextCell 'D20' '0.0₃0827'
Set-TextCell 'E20' '  -3.77%  '
Set-TextCell 'D21' '233.63'
Set-TextCell 'E21' '  +1.61%  '
Set-TextCell 'D22' '5.05'
Set-TextCell 'E22' '  -2.88%  '
Set-TextCell 'E23' '  -0.05%  '
Set-TextCell 'D24' '2.54'
Set-TextCell 'E24' '  -0.31%  '
Set-TextCell 'D25' '2.37'
Set-TextCell 'E25' '  +0.99%  '
Set-TextCell 'D26' '164.41'
Set-TextCell 'E26' '  +0.73%  '
Set-TextCell 'D27' '8.87'
Set-TextCell 'E27' '  -5.05%  '
Set-TextCell 'D28' '19.43'
Set-TextCell 'E28' '  -1.35%  '
Set-TextCell 'B29' 'Kaspa'
Set-TextCell 'C29' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D29' '0.125'
Set-TextCell 'E29' '  -8.55%  '
Set-TextCell 'B30' 'ImmutableX'
Set-TextCell 'C30' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D30' '1.34'
Set-TextCell 'E30' '  -2.94%  '
Set-TextCell 'E31' '  -2.03%  '
Set-TextCell 'D32' '4.56'
Set-TextCell 'E32' '  -3.99%  '
Set-TextCell 'D33' '0.0636'
Set-TextCell 'E33' '  -5.13%  '
Set-TextCell 'D34' '4.38'
Set-TextCell 'E34' '  -1.43%  '
Set-TextCell 'E35' '  -6.45%  '
Set-TextCell 'D36' '3.49'
Set-TextCell 'E36' '  -2.70%  '
Set-TextCell 'E37' '  +0.91%  '
Set-TextCell 'E38' '  +0.03%  '
Set-TextCell 'D39' '5.48'
Set-TextCell 'E39' '  +2.99%  '
Set-TextCell 'E40' '  +1.32%  '
Set-TextCell 'E41' '  -0.40%  '
Set-TextCell 'D42' '1.438.07'
Set-TextCell 'E42' '  +4.88%  '
Set-TextCell 'B43' 'VeChain'
Set-TextCell 'C43' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D43' '0.0209'
Set-TextCell 'E43' '  -3.22%  '
Set-TextCell 'B44' 'Cronos'
Set-TextCell 'C44' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D44' '0.0912'
Set-TextCell 'E44' '  -5.93%  '
Set-TextCell 'D45' '88.77'
Set-TextCell 'E45' '  -2.50%  '
Set-TextCell 'D46' '15.47'
Set-TextCell 'E46' '  -6.83%  '
Set-TextCell 'E47' '  -3.30%  '
Set-TextCell 'E48' '  +1.53%  '
Set-TextCell 'D49' '6.88'
Set-TextCell 'E49' '  -7.21%  '
Set-TextCell 'D50' '2.190.61'
Set-TextCell 'E50' '  -0.71%  '
Set-TextCell 'E51' '  -8.37%  '
